$d = $word.ActiveDocument

# --- Step 1: remove the trailing empty paragraph ---------------------------
# The document has two paragraphs: the content paragraph ("requirements.txt")
# and a second, completely empty one. Word removes a paragraph by deleting
# its preceding paragraph mark, which merges it into the paragraph before it.
if ($d.Paragraphs.Count -gt 1) {
    $first = $d.Paragraphs.Item(1)
    $markPos = $first.Range.End
    $mergeRange = $d.Range($markPos - 1, $markPos + 1)
    $mergeRange.Delete()
}

# --- Step 2: replace the remaining paragraph's content ----------------------
# Swap the formatted "requirements.txt" run (and the paragraph-level shading/
# tabs/spacing/font formatting that decorated it) for a single plain run
# containing the new URL, keeping the existing _GoBack bookmark - but now
# positioned after the run instead of before it.
$p1 = $d.Paragraphs.Item(1)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>https://github.com/streamlit/streamlit.git</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
[void]$p1.Range.InsertXML($newParaXml)
